$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header label rename: "Start" -> "Zones"
$ws.Range("A1").Value = "Zones"

# Make sure any newly-populated cells in the matrix body pick up the same
# formatting (style) already used by the existing matrix cells (e.g. C2),
# then fill in the previously-empty lower-triangle/diagonal cells so the
# table becomes a full symmetric zone distance/time matrix (diagonal = 1,
# i.e. same zone).
$ws.Range("C2").Copy() | Out-Null
$ws.Range("B2:J10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$values = @{
    "B2" = 1;    "C2" = 1.5;  "D2" = 2;    "E2" = 2;    "F2" = 2;    "G2" = 3;    "H2" = 4;    "I2" = 4;    "J2" = 3;
    "B3" = 1.5;  "C3" = 1;    "D3" = 1.5;  "E3" = 2;    "F3" = 2;    "G3" = 3;    "H3" = 4;    "I3" = 4;    "J3" = 3;
    "B4" = 2;    "C4" = 1.5;  "D4" = 1;    "E4" = 1.5;  "F4" = 2;    "G4" = 3;    "H4" = 3;    "I4" = 4;    "J4" = 3;
    "B5" = 2;    "C5" = 2;    "D5" = 1.5;  "E5" = 1;    "F5" = 1.5;  "G5" = 2;    "H5" = 2;    "I5" = 3;    "J5" = 2;
    "B6" = 2;    "C6" = 2;    "D6" = 2;    "E6" = 1.5;  "F6" = 1;    "G6" = 2;    "H6" = 2;    "I6" = 3;    "J6" = 2;
    "B7" = 3;    "C7" = 3;    "D7" = 3;    "E7" = 2;    "F7" = 2;    "G7" = 1;    "H7" = 1.5;  "I7" = 2;    "J7" = 2.5;
    "B8" = 4;    "C8" = 4;    "D8" = 3;    "E8" = 2;    "F8" = 2;    "G8" = 1.5;  "H8" = 1;    "I8" = 1.5;  "J8" = 2;
    "B9" = 4;    "C9" = 4;    "D9" = 4;    "E9" = 3;    "F9" = 3;    "G9" = 2;    "H9" = 1.5;  "I9" = 1;    "J9" = 3;
    "B10" = 3;   "C10" = 3;   "D10" = 3;   "E10" = 2;   "F10" = 2;   "G10" = 2.5; "H10" = 2;   "I10" = 3;   "J10" = 1;
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
